$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Value = -7.995
$ws.Range("D18").Value = -8.377000000000001
$ws.Range("D20").Value = -7.619999999999999
$ws.Range("D27").Value = -7.866000000000001
$ws.Range("D69").Value = -7.678999999999999
$ws.Range("D76").Value = -7.687
$ws.Range("D82").Value = -8.17
